# Update imputed values in result_data_KNN sheet ("Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.924
$ws.Range("D10").Value = -7.627000000000001
$ws.Range("D15").Value = -7.872
$ws.Range("C18").Value = -12.525
$ws.Range("A21").Value = -20.091
$ws.Range("D21").Value = -8.297999999999998
$ws.Range("D22").Value = -7.870000000000002
$ws.Range("A23").Value = -20.808
$ws.Range("B24").Value = 5.572000000000001
$ws.Range("D24").Value = -7.581
$ws.Range("A25").Value = -21.63
$ws.Range("B28").Value = 5.548999999999999
$ws.Range("B36").Value = 7.431000000000002
$ws.Range("B45").Value = 5.655999999999999
$ws.Range("D46").Value = -7.969999999999999
$ws.Range("B48").Value = 5.48
$ws.Range("B49").Value = 6.269
$ws.Range("C51").Value = -11.589
$ws.Range("B52").Value = 5.280000000000001
$ws.Range("A53").Value = -20.888
$ws.Range("B53").Value = 6.704000000000001
$ws.Range("B54").Value = 5.197
$ws.Range("C55").Value = -13.639
$ws.Range("D56").Value = -7.607000000000001
$ws.Range("A57").Value = -22.178
$ws.Range("A59").Value = -22.461
$ws.Range("D61").Value = -8.052000000000001
$ws.Range("C64").Value = -10.76
$ws.Range("D66").Value = -7.290000000000001
$ws.Range("A69").Value = -21.53
$ws.Range("B70").Value = 4.935
$ws.Range("D74").Value = -7.958999999999999
$ws.Range("A79").Value = -21.186
$ws.Range("C80").Value = -12.033
$ws.Range("A83").Value = -22.08
$ws.Range("B86").Value = 5.197
$ws.Range("B87").Value = 4.55
$ws.Range("D87").Value = -8.004000000000001
$ws.Range("D88").Value = -7.597
$ws.Range("C92").Value = -11.072
$ws.Range("A93").Value = -21.598
$ws.Range("C94").Value = -11.945
$ws.Range("C96").Value = -11.199
$ws.Range("D100").Value = -7.295000000000002
$ws.Range("B101").Value = 5.197
